$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.241.20"
$ws.Range("E2").Value = "  +2.01%  "
$ws.Range("D3").Value = "2.360.08"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "542.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.71%  "
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.563"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.99%  "
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("E10").Value = "  +5.75%  "
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("E12").Value = "  +2.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").Value = "2.780.57"
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("D15").Value = "58.217.94"
$ws.Range("E15").Value = "  +1.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000134"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").Value = "2.357.72"
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "332.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "62.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.16%  "
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("E25").Value = "  -2.63%  "
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("E27").Value = "  +5.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  +2.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("E33").Value = "  +12.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.28"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.33%  "
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "39.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "148.09"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "294.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.378"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0950"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0504"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("E47").Value = "  +1.60%  "
$ws.Range("E48").Value = "  +1.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.382"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.54%  "
